$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A4").Value = -21.18349999999998
$ws.Range("A6").Value = -22.70770000000001
$ws.Range("A7").Value = -20.41029999999999
$ws.Range("D7").Value = -7.406500000000003
$ws.Range("D12").Value = -7.181900000000002
$ws.Range("E13").Value = 16.747
$ws.Range("E14").Value = 17.08290000000001
$ws.Range("D15").Value = -8.361099999999997
$ws.Range("A16").Value = -21.81380000000001
$ws.Range("E16").Value = 16.04340000000001
$ws.Range("E19").Value = 16.3752
$ws.Range("A20").Value = -20.49099999999999
$ws.Range("D20").Value = -7.771899999999995
$ws.Range("D21").Value = -7.901599999999994
$ws.Range("D22").Value = -7.994900000000005
$ws.Range("E22").Value = 16.3159
$ws.Range("D23").Value = -7.191999999999994
$ws.Range("A28").Value = -22.0165
$ws.Range("A29").Value = -21.36969999999998
$ws.Range("D29").Value = -7.2155
$ws.Range("A32").Value = -21.08889999999998
$ws.Range("D34").Value = -7.963899999999998
$ws.Range("E36").Value = 15.79900000000001
$ws.Range("A40").Value = -20.42850000000002
$ws.Range("D42").Value = -8.507600000000005
$ws.Range("D43").Value = -8.174300000000001
$ws.Range("D44").Value = -7.597199999999997
$ws.Range("D45").Value = -7.550699999999999
$ws.Range("A46").Value = -21.90110000000002
$ws.Range("D46").Value = -7.638199999999992
$ws.Range("E46").Value = 17.11920000000001
$ws.Range("D50").Value = -8.321899999999994
$ws.Range("E50").Value = 16.6662
$ws.Range("A51").Value = -21.9949
$ws.Range("D51").Value = -7.529799999999999
$ws.Range("A52").Value = -22.04699999999998
$ws.Range("A57").Value = -22.62320000000001
$ws.Range("A59").Value = -21.92509999999999
$ws.Range("A62").Value = -22.0657
$ws.Range("A66").Value = -21.36079999999999
$ws.Range("D66").Value = -7.244200000000002
$ws.Range("D67").Value = -6.530800000000004
$ws.Range("A73").Value = -20.06959999999999
$ws.Range("A74").Value = -21.92409999999999
$ws.Range("D79").Value = -6.379000000000007
$ws.Range("D84").Value = -8.844100000000005
$ws.Range("A92").Value = -21.5599
$ws.Range("D92").Value = -6.474300000000004
$ws.Range("E95").Value = 18.03270000000002
$ws.Range("D97").Value = -8.042699999999996
$ws.Range("E97").Value = 16.88259999999999
$ws.Range("A100").Value = -22.0322
